$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("H2").Phonetics.Font.Size = 8
$ws.Range("H2").Font.Size = 12
